$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 2 (shifts existing data rows 2-18 down to 3-19),
# then strip any inherited formatting so the new row matches the plain
# (unstyled) data rows used throughout the rest of the sheet.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# The A/B/C columns hold dash-separated dates stored as plain text
# (e.g. "2024-06-04"). Force text format first so Excel doesn't
# auto-convert them into date serial numbers, then write the values.
$dateRange = $ws.Range("A2:C2")
$dateRange.NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2024-06-04"
$ws.Cells.Item(2, 2).Value = "2024-06-05"
$ws.Cells.Item(2, 3).Value = "2024-06-21"
$dateRange.ClearFormats()

# Remaining fields of the new row (new IPO record: KB제29호스팩).
$ws.Cells.Item(2, 4).Value = "KB"
$ws.Cells.Item(2, 5).Value = "KB제29호스팩"
$ws.Cells.Item(2, 6).Value = 6000000
$ws.Cells.Item(2, 7).Value = 6000000
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 2000
$ws.Cells.Item(2, 10).Value = 2000
$ws.Cells.Item(2, 11).Value = 6220000
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 2000
$ws.Cells.Item(2, 14).Value = "1100.48:1"
$ws.Cells.Item(2, 15).Value = "-"
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = "기업인수합병"

# The old last row (originally row 18, 신한제13호스팩) was pushed down to row 19
# by the insert above; the dataset drops that record entirely, so remove it.
$ws.Rows.Item(19).Delete()
